$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New defined names (so formulas can refer to ranges by name) ---
$wb.Names.Add('Price', '=Sheet1!$H$6:$H$9')
$wb.Names.Add('Prices', '=Sheet1!$H$6:$H$9')
$wb.Names.Add('Quantity', '=Sheet1!$I$6:$I$9')
$wb.Names.Add('Cost', '=Sheet1!$J$6:$J$9')

# --- Column width for the new Cost column ---
$ws.Columns.Item(10).ColumnWidth = 13

# --- Headers (row 5): Price / Quantity / Cost ---
$ws.Range('H5').Value = 'Price'
$ws.Range('I5').Value = 'Quantity '
$ws.Range('J5').Value = 'Cost'
$ws.Range('H5:J5').Style = 'Heading 3'
$ws.Range('H5:J5').HorizontalAlignment = -4108

# --- Price data (H6:H9) ---
$ws.Range('H6').Value = 2.99
$ws.Range('H7').Value = 34.99
$ws.Range('H8').Value = 42.5
$ws.Range('H9').Value = 56.13

# --- Quantity data (I6:I9) ---
$ws.Range('I6').Value = 5
$ws.Range('I7').Value = 65
$ws.Range('I8').Value = 45
$ws.Range('I9').Value = 35

# --- Cost = Price * Quantity (dynamic array formula, spills J6:J9) ---
$ws.Range('J6').Formula2 = '=Price*Quantity'

# --- Number formats matching the original author's formatting ---
$ws.Range('H6:H8').NumberFormat = '_("$"* #,##0.00_);_("$"* \(#,##0.00\);_("$"* "-"??_);_(@_)'
$ws.Range('J6:J8').NumberFormat = '_("$"* #,##0.00_);_("$"* \(#,##0.00\);_("$"* "-"??_);_(@_)'
$ws.Range('H9').NumberFormat = '_-[$$-409]* #,##0.00_ ;_-[$$-409]* \-#,##0.00\ ;_-[$$-409]* "-"??_ ;_-@_ '
$ws.Range('J9').Style = 'Currency'

$ws.Range('I6:I9').HorizontalAlignment = -4108

# --- Selection state, matching the saved workbook ---
$ws.Range('N12').Select() | Out-Null
